$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 182.4
$ws.Range("I33").Value = 182.4
$ws.Range("K33").Value = 182.4
$ws.Range("M33").Value = 46.59999999999999
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H94").Value = 10000
$ws.Range("J94").Value = 7866.6665
$ws.Range("L94").Value = 7866.6665
$ws.Range("N94").Value = -8768.666499999999
$ws.Range("H127").Value = 9163.333000000001
$ws.Range("I127").Value = 9163.333000000001
$ws.Range("K127").Value = 27489.999
$ws.Range("M127").Value = -22529.999
$ws.Range("H137").Value = 1947.3
$ws.Range("I137").Value = 1912.1666
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 5736.4998
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = -3186.4998
$ws.Range("N137").Value = -11100
$ws.Range("H138").Value = 4694.0166
$ws.Range("I138").Value = 3273.7646
$ws.Range("J138").Value = 6482.4814
$ws.Range("K138").Value = 9821.293799999999
$ws.Range("L138").Value = 19447.4442
$ws.Range("M138").Value = -4681.293799999999
$ws.Range("N138").Value = -29727.4442
$ws.Range("H141").Value = 2055.3438
$ws.Range("I141").Value = 1805.8966
$ws.Range("K141").Value = 5417.6898
$ws.Range("M141").Value = -237.6898000000001

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 159.5
$ws.Range("I4").Value = 119
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 119
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = -3
$ws.Range("N4").Value = -432
$ws.Range("H5").Value = 76
$ws.Range("I5").Value = 50
$ws.Range("K5").Value = 50
$ws.Range("M5").Value = 62
$ws.Range("H114").Value = 26398
$ws.Range("J114").Value = 26398
$ws.Range("L114").Value = 26398
$ws.Range("N114").Value = -35076
$ws.Range("H139").Value = 90000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 90000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 90000
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -100280
$ws.Range("H140").Value = 95500
$ws.Range("J140").Value = 95500
$ws.Range("L140").Value = 95500
$ws.Range("N140").Value = -105860

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 76
$ws.Range("I4").Value = 50
$ws.Range("K4").Value = 50
$ws.Range("M4").Value = 65
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 1643.2307
$ws.Range("I134").Value = 1643.2307
$ws.Range("K134").Value = 4929.6921
$ws.Range("M134").Value = -2394.6921

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2362.05
$ws.Range("I58").Value = 1298.0625
$ws.Range("J58").Value = 6618
$ws.Range("K58").Value = 1298.0625
$ws.Range("L58").Value = 6618
$ws.Range("M58").Value = -1095.0625
$ws.Range("N58").Value = -7024
$ws.Range("H60").Value = 13118.125
$ws.Range("I60").Value = 10659.4
$ws.Range("K60").Value = 10659.4
$ws.Range("M60").Value = -10148.4
$ws.Range("H132").Value = 1067.9
$ws.Range("I132").Value = 1094.2858
$ws.Range("J132").Value = 1006.3333
$ws.Range("K132").Value = 3282.8574
$ws.Range("L132").Value = 3018.9999
$ws.Range("M132").Value = -752.8574000000003
$ws.Range("N132").Value = -8078.9999
$ws.Range("H136").Value = 2362.05
$ws.Range("I136").Value = 1298.0625
$ws.Range("J136").Value = 6618
$ws.Range("K136").Value = 3894.1875
$ws.Range("L136").Value = 19854
$ws.Range("M136").Value = -1344.1875
$ws.Range("N136").Value = -24954

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1753159.1
$ws.Range("I4").Value = 2561903.5
$ws.Range("K4").Value = 7685710.5
$ws.Range("M4").Value = -7685598.5
$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 90000
$ws.Range("N123").Value = -94900

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 11998
$ws.Range("I33").Value = 11994
$ws.Range("J33").Value = 11999.333
$ws.Range("K33").Value = 11994
$ws.Range("L33").Value = 11999.333
$ws.Range("M33").Value = -11742
$ws.Range("N33").Value = -12503.333
$ws.Range("H44").Value = 13967
$ws.Range("J44").Value = 13967
$ws.Range("L44").Value = 13967
$ws.Range("N44").Value = -15159
$ws.Range("H109").Value = 5316.6665
$ws.Range("J109").Value = 5316.6665
$ws.Range("L109").Value = 5316.6665
$ws.Range("N109").Value = -7396.6665
$ws.Range("H134").Value = 63662.25
$ws.Range("J134").Value = 63662.25
$ws.Range("L134").Value = 190986.75
$ws.Range("N134").Value = -196056.75

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2363.5
$ws.Range("J7").Value = 1700
$ws.Range("L7").Value = 1700
$ws.Range("N7").Value = -1924
$ws.Range("H40").Value = 2574.0908
$ws.Range("I40").Value = 2331.5
$ws.Range("K40").Value = 2331.5
$ws.Range("M40").Value = -2195.5
$ws.Range("H61").Value = 6401
$ws.Range("I61").Value = 7998.6665
$ws.Range("J61").Value = 4004.5
$ws.Range("K61").Value = 7998.6665
$ws.Range("L61").Value = 4004.5
$ws.Range("M61").Value = -7796.6665
$ws.Range("N61").Value = -4408.5
$ws.Range("H113").Value = 6401
$ws.Range("I113").Value = 7998.6665
$ws.Range("J113").Value = 4004.5
$ws.Range("K113").Value = 7998.6665
$ws.Range("L113").Value = 4004.5
$ws.Range("M113").Value = -5828.6665
$ws.Range("N113").Value = -8344.5
$ws.Range("H122").Value = 4155.75
$ws.Range("I122").Value = 3541.1667
$ws.Range("K122").Value = 10623.5001
$ws.Range("M122").Value = -8173.500100000001
$ws.Range("H126").Value = 2363.5
$ws.Range("J126").Value = 1700
$ws.Range("L126").Value = 5100
$ws.Range("N126").Value = -10040

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 7500.4443
$ws.Range("I14").Value = 2584
$ws.Range("J14").Value = 17333.334
$ws.Range("K14").Value = 2584
$ws.Range("L14").Value = 17333.334
$ws.Range("M14").Value = -2416
$ws.Range("N14").Value = -17669.334
$ws.Range("H52").Value = 36499.5
$ws.Range("I52").Value = 36499.5
$ws.Range("K52").Value = 36499.5
$ws.Range("M52").Value = -36273.5
$ws.Range("H80").Value = 65000
$ws.Range("J80").Value = 65000
$ws.Range("L80").Value = 65000
$ws.Range("N80").Value = -66996
$ws.Range("H83").Value = 65000
$ws.Range("J83").Value = 65000
$ws.Range("L83").Value = 195000
$ws.Range("N83").Value = -204984
$ws.Range("H93").Value = 50000
$ws.Range("J93").Value = 50000
$ws.Range("L93").Value = 50000
$ws.Range("N93").Value = -54992
$ws.Range("H126").Value = 3068.8
$ws.Range("I126").Value = 2872.818
$ws.Range("J126").Value = 3308.3333
$ws.Range("K126").Value = 8618.454000000002
$ws.Range("L126").Value = 9924.999899999999
$ws.Range("M126").Value = -6148.454000000002
$ws.Range("N126").Value = -14864.9999
$ws.Range("H132").Value = 20884.277
$ws.Range("I132").Value = 30763.516
$ws.Range("J132").Value = 1743.25
$ws.Range("K132").Value = 92290.548
$ws.Range("L132").Value = 5229.75
$ws.Range("M132").Value = -89760.548
$ws.Range("N132").Value = -10289.75
$ws.Range("H136").Value = 1215.6923
$ws.Range("I136").Value = 1232.32
$ws.Range("K136").Value = 3696.96
$ws.Range("M136").Value = -1146.96
